$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for MOUDERY (row 4) and DINDEFELO (row 9), matching format of surrounding rows
$ws.Rows("4:4").Insert()
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A4").Borders.LineStyle = 1

$ws.Rows("9:9").Insert()
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4160
$ws.Range("A9").Borders.LineStyle = 1

# Write the final data for every row (labels + counts + percentages)
$ws.Range("A2").Value = "BALLOU"
$ws.Range("B2").Value = 450
$ws.Range("C2").Value = 5060
$ws.Range("D2").Value = 5510
$ws.Range("E2").Value = 8.1669691470054442
$ws.Range("F2").Value = 91.833030852994554

$ws.Range("A3").Value = "GABOU"
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = 3815
$ws.Range("D3").Value = 4015
$ws.Range("E3").Value = 4.9813200498132
$ws.Range("F3").Value = 95.018679950186808

$ws.Range("A4").Value = "MOUDERY"
$ws.Range("B4").Value = 49
$ws.Range("C4").Value = 953
$ws.Range("D4").Value = 1002
$ws.Range("E4").Value = 4.8902195608782426
$ws.Range("F4").Value = 95.109780439121764

$ws.Range("A5").Value = "NDOGA BABACAR"
$ws.Range("B5").Value = 97
$ws.Range("C5").Value = 1081
$ws.Range("D5").Value = 1178
$ws.Range("E5").Value = 8.2342954159592523
$ws.Range("F5").Value = 91.765704584040748

$ws.Range("A6").Value = "MISSIRAH"
$ws.Range("B6").Value = 1015
$ws.Range("C6").Value = 5076
$ws.Range("D6").Value = 6091
$ws.Range("E6").Value = 16.66393038909867
$ws.Range("F6").Value = 83.336069610901333

$ws.Range("A7").Value = "NETTEBOULOU"
$ws.Range("B7").Value = 484
$ws.Range("C7").Value = 6427
$ws.Range("D7").Value = 6911
$ws.Range("E7").Value = 7.0033280277817962
$ws.Range("F7").Value = 92.996671972218209

$ws.Range("A8").Value = "BANDAFASSI"
$ws.Range("B8").Value = 2381
$ws.Range("C8").Value = 5784
$ws.Range("D8").Value = 8165
$ws.Range("E8").Value = 29.161053276178809
$ws.Range("F8").Value = 70.838946723821181

$ws.Range("A9").Value = "DINDEFELO"
$ws.Range("B9").Value = 363
$ws.Range("C9").Value = 1167
$ws.Range("D9").Value = 1530
$ws.Range("E9").Value = 23.725490196078429
$ws.Range("F9").Value = 76.274509803921561

$ws.Range("A10").Value = "DIMBOLI"
$ws.Range("B10").Value = 1353
$ws.Range("C10").Value = 3504
$ws.Range("D10").Value = 4857
$ws.Range("E10").Value = 27.856701667696111
$ws.Range("F10").Value = 72.1432983323039

$ws.Range("A11").Value = "FONGOLIMBI"
$ws.Range("B11").Value = 1344
$ws.Range("C11").Value = 3041
$ws.Range("D11").Value = 4385
$ws.Range("E11").Value = 30.649942987457241
$ws.Range("F11").Value = 69.350057012542749

$ws.Range("A12").Value = "BEMBOU"
$ws.Range("B12").Value = 800
$ws.Range("C12").Value = 2650
$ws.Range("D12").Value = 3450
$ws.Range("E12").Value = 23.188405797101449
$ws.Range("F12").Value = 76.811594202898547

# Match the final selected cell
$ws.Range("G11").Select() | Out-Null
